$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, shifting existing rows 43:103 down to 44:104.
$ws.Rows("43:43").Insert()

# Populate the newly inserted row 43. This is a new weekly price record for the
# same market/product series; it carries the same reference data (market,
# region, product, quality, price range, unit, origin, weight factor and
# classification) as the neighboring historical entries, with a new date (D)
# and a new volume (J).
$ws.Range("A43").Value = 11
$ws.Range("B43").Value = "Vega Monumental Concepción"
$ws.Range("C43").Value = "Bíobío"
$ws.Range("D43").Value = 44966
$ws.Range("E43").Value = 8
$ws.Range("F43").Value = 100112012
$ws.Range("G43").Value = "Espinaca"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 40
$ws.Range("K43").Value = 7000
$ws.Range("L43").Value = 7500
$ws.Range("M43").Value = 7250
$ws.Range("N43").Value = "$/cuna 10 kilos"
$ws.Range("O43").Value = "Región Metropolitana"
$ws.Range("P43").Value = 725
$ws.Range("Q43").Value = 10
$ws.Range("R43").Value = "Hortaliza"
